$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.671.52"
$ws.Range("E2").Value = "  -0.49%  "
$ws.Range("D3").Value = "1.899.82"
$ws.Range("E3").Value = "  -0.05%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9992"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.86%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.61"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9986"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.83%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5023"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +4.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3764"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07235"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.94"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.52%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8911"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.78%  "
$ws.Range("D12").Value = "1.944.96"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07621"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.434"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.73%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.61"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9997"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008757"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.91%  "
$ws.Range("E18").Value = "  -0.69%  "
$ws.Range("D19").Value = "27.711.70"
$ws.Range("E19").Value = "  -0.62%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").Value = "2.136.63"
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("E23").Value = "  -0.35%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.574"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.46%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "153.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.847"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.67%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.176"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.821"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08930"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.184"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.229"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.05%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7835"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.79%  "
$ws.Range("B35").Value = "Filecoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.783"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.621"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +4.29%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02082"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.14%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.054"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.00%  "
$ws.Range("E39").Value = "  -0.66%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.5484"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.58%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.05274"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.06%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.704"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.40%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.75%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.431"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.13%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1506"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4765"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.19%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9985"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.87%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.610"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.87%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.48"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.65%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05996"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.45%  "
